$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<all>"
$ws.Range("C2").Value = 19

# Row 3
$ws.Range("C3").Value = 20

# Row 4
$ws.Range("C4").Value = 19

# Row 5
$ws.Range("C5").Value = 20

# Row 6
$ws.Range("B6").Value = "<sen>"
$ws.Range("C6").Value = 15

# Row 7
$ws.Range("C7").Value = 10

# Row 8
$ws.Range("C8").Value = 19

# Row 9
$ws.Range("C9").Value = 13

# Row 10
$ws.Range("C10").Value = 13

# Row 11
$ws.Range("C11").Value = 20

# Row 13
$ws.Range("C13").Value = 19

# Row 14
$ws.Range("B14").Value = "<all>"
$ws.Range("C14").Value = 12

# Row 15
$ws.Range("C15").Value = 9

# Row 17
$ws.Range("C17").Value = 17

# Row 18
$ws.Range("C18").Value = 13
